$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1: same header-row formatting (bold, centered, bordered)
# as the existing G1 header -- copy its format only, then set the text.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# New data column values
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
